$d = $word.ActiveDocument

# ---- Table 1 (Отделение ИТБ) ----
$tbl1 = $d.Tables.Item(1)

# Group header "ИС-221"
$cell = $tbl1.Rows.Item(3).Cells.Item(2)
$cell.Range.Text = "ИС-221"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 14
$cell.Range.Font.SizeBi = 14
$cell.Range.Font.Underline = 1

# "3п. – нет " note for ИС-221
$cell = $tbl1.Rows.Item(4).Cells.Item(2)
$cell.Range.Text = "3п. – нет "
$cell.Range.Font.Italic = 1

# ---- Table 2 (Отделение СРПП) ----
$tbl2 = $d.Tables.Item(2)

# Group header "СР-231"
$cell = $tbl2.Rows.Item(3).Cells.Item(2)
$cell.Range.Text = "СР-231"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 14
$cell.Range.Font.SizeBi = 14
$cell.Range.Font.Underline = 1

# Group header "ДО-211"
$cell = $tbl2.Rows.Item(3).Cells.Item(3)
$cell.Range.Text = "ДО-211"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 14
$cell.Range.Font.SizeBi = 14
$cell.Range.Font.Underline = 1

# "Дистанционные занятия" note for СР-231
$cell = $tbl2.Rows.Item(4).Cells.Item(2)
$cell.Range.Text = "Дистанционные занятия"
$cell.Range.Font.Italic = 1

# ДО-211 schedule note (two paragraphs, one with an in-line break)
$cell = $tbl2.Rows.Item(4).Cells.Item(3)
$cell.Range.Text = "4п. Дет. Лит-ра с ПВЧ" + [char]13 + "Цыпдакова, 512 ауд" + [char]11 + "6п. – нет "
$cell.Range.Font.Italic = 1

# Group header "ДО-212"
$cell = $tbl2.Rows.Item(5).Cells.Item(1)
$cell.Range.Text = "ДО-212"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 14
$cell.Range.Font.SizeBi = 14
$cell.Range.Font.Underline = 1

# ДО-212 schedule note (two paragraphs)
$cell = $tbl2.Rows.Item(6).Cells.Item(1)
$cell.Range.Text = "5п. Дет. Лит-ра с ПВЧ" + [char]13 + "Цыпдакова, 512 ауд"
$cell.Range.Font.Italic = 1

Write-Output "edit complete"
